$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment
$ws.Range('D2').Value = '40.118.67'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').Value = '2.346.82'
$ws.Range('E3').Value = '  -3.14%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  -1.76%  '
$ws.Range('E6').Value = '  -3.77%  '
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -2.22%  '
$ws.Range('E10').Value = '  -2.20%  '
$ws.Range('E11').Value = '  -5.78%  '
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('D13').Value = '2.709.48'
$ws.Range('E13').Value = '  -3.11%  '
$ws.Range('E14').Value = '  -3.81%  '
$ws.Range('E15').Value = '  -5.50%  '
$ws.Range('D16').Value = '2.388.85'
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '40.099.60'
$ws.Range('E18').Value = '  -2.12%  '
$ws.Range('E19').Value = '  -1.78%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('E21').Value = '  -5.06%  '
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  -4.82%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  -2.39%  '
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('E36').Value = '  -0.71%  '
$ws.Range('E37').Value = '  -3.60%  '
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('E39').Value = '  -6.35%  '
$ws.Range('E40').Value = '  -2.38%  '
$ws.Range('E41').Value = '  +0.48%  '
$ws.Range('D42').Value = '1.965.12'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('E43').Value = '  -2.91%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E44').Value = '  -4.41%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E45').Value = '  -3.21%  '
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  -5.82%  '
$ws.Range('D48').Value = '2.569.75'
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('E51').Value = '  -2.35%  '

# Numeric-looking text values: must preserve as text (leading/trailing zeros, dot format)
# Temporarily force text format, assign, then restore original style to avoid residual formatting diffs
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '310.51'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '85.51'
$ws.Range('D6').Style = $origStyle
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.530'
$ws.Range('D7').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0810'
$ws.Range('D10').Style = $origStyle
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '30.03'
$ws.Range('D11').Style = $origStyle
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.43'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.78'
$ws.Range('D15').Style = $origStyle
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.759'
$ws.Range('D17').Style = $origStyle
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.10'
$ws.Range('D20').Style = $origStyle
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '68.08'
$ws.Range('D21').Style = $origStyle
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '10.71'
$ws.Range('D22').Style = $origStyle
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.39'
$ws.Range('D23').Style = $origStyle
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.71'
$ws.Range('D27').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.25'
$ws.Range('D29').Style = $origStyle
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '34.57'
$ws.Range('D30').Style = $origStyle
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '153.74'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.11'
$ws.Range('D33').Style = $origStyle
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0720'
$ws.Range('D35').Style = $origStyle
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.82'
$ws.Range('D37').Style = $origStyle
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0983'
$ws.Range('D38').Style = $origStyle
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '15.57'
$ws.Range('D39').Style = $origStyle
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.87'
$ws.Range('D41').Style = $origStyle
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '17.75'
$ws.Range('D44').Style = $origStyle
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0266'
$ws.Range('D45').Style = $origStyle
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '9.36'
$ws.Range('D46').Style = $origStyle
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.70'
$ws.Range('D47').Style = $origStyle
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '93.12'
$ws.Range('D49').Style = $origStyle
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '70.53'
$ws.Range('D50').Style = $origStyle
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '50.46'
$ws.Range('D51').Style = $origStyle
